$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" '65.752.64'
Set-TextValue "E2" '  -1.25%  '
Set-TextValue "D3" '3.418.34'
Set-TextValue "E3" '  -1.50%  '
Set-TextValue "E4" '  -0.04%  '
Set-TextValue "D5" '596.49'
Set-TextValue "E5" '  -0.57%  '
Set-TextValue "D6" '142.32'
Set-TextValue "E6" '  -3.14%  '
Set-TextValue "D7" '3.417.23'
Set-TextValue "E7" '  -1.48%  '
Set-TextValue "E8" '  -0.26%  '
Set-TextValue "E9" '  -1.97%  '
Set-TextValue "D10" '8.06'
Set-TextValue "E10" '  +6.74%  '
Set-TextValue "E11" '  -5.32%  '
Set-TextValue "E12" '  -3.91%  '
Set-TextValue "D13" '3.995.97'
Set-TextValue "E13" '  -1.47%  '
Set-TextValue "D14" '0.0000200'
Set-TextValue "E14" '  -5.68%  '
Set-TextValue "D15" '29.61'
Set-TextValue "E15" '  -5.40%  '
Set-TextValue "D16" '3.417.22'
Set-TextValue "E16" '  -1.55%  '
Set-TextValue "E17" '  -0.68%  '
Set-TextValue "D18" '65.816.52'
Set-TextValue "E18" '  -1.30%  '
Set-TextValue "D19" '10.30'
Set-TextValue "E19" '  +3.03%  '
Set-TextValue "D20" '6.11'
Set-TextValue "E20" '  -4.53%  '
Set-TextValue "D21" '14.56'
Set-TextValue "E21" '  -4.71%  '
Set-TextValue "D22" '415.28'
Set-TextValue "E22" '  -4.27%  '
Set-TextValue "E23" '  -4.85%  '
Set-TextValue "D24" '77.22'
Set-TextValue "E25" '  +0.16%  '
Set-TextValue "E26" '  -8.38%  '
Set-TextValue "E27" '  -4.92%  '
Set-TextValue "D28" '7.89'
Set-TextValue "E28" '  -5.45%  '
Set-TextValue "E29" '  -2.20%  '
Set-TextValue "E31" '  -4.15%  '
Set-TextValue "E32" '  -8.10%  '
Set-TextValue "D33" '24.58'
Set-TextValue "E33" '  -2.76%  '
Set-TextValue "D34" '3.414.54'
Set-TextValue "E34" '  -1.37%  '
Set-TextValue "E36" '  -6.58%  '
Set-TextValue "D37" '5.47'
Set-TextValue "E37" '  -7.93%  '
Set-TextValue "E38" '  -4.59%  '
Set-TextValue "E39" '  -0.02%  '
Set-TextValue "D40" '168.60'
Set-TextValue "E40" '  -4.24%  '
Set-TextValue "E41" '  -2.98%  '
Set-TextValue "D42" '0.873'
Set-TextValue "E42" '  -2.01%  '
Set-TextValue "D43" '5.03'
Set-TextValue "E43" '  -6.90%  '
Set-TextValue "E44" '  -10.40%  '
Set-TextValue "D45" '45.37'
Set-TextValue "E45" '  -2.07%  '
Set-TextValue "D46" '26.32'
Set-TextValue "E46" '  -8.64%  '
Set-TextValue "D47" '1.18'
Set-TextValue "E47" '  -3.78%  '
Set-TextValue "D48" '7.05'
Set-TextValue "E48" '  -4.98%  '
Set-TextValue "D49" '2.26'
Set-TextValue "E49" '  -6.57%  '
Set-TextValue "E50" '  -5.99%  '
Set-TextValue "D51" '0.231'
Set-TextValue "E51" '  -5.55%  '
